{"js": "// Add a new \"Dev Ops - Moderate\" bullet to the Technical-Skills list, right\n// before the \"OOMD (Object Oriented Modeling and Design) - Basics\" item.\n// \"Dev Ops\" is regular weight, \" - Moderate\" is bold - matching the style\n// already used by every other skill line in this list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its (stable) text rather than a hard-coded\n// index so the script keeps working even if earlier content shifts.\nconst anchor = paragraphs.items.find((p) =>\n  p.text.indexOf(\"OOMD (Object Oriented Modeling and Design)\") !== -1\n);\nif (!anchor) {\n  throw new Error('Could not find the \"OOMD (Object Oriented Modeling and Design)\" paragraph.');\n}\n\n// Insert a brand-new paragraph right before it; Word copies the anchor's\n// paragraph formatting (numbering, indent, fonts, size) automatically.\nconst newPara = anchor.insertParagraph(\"Dev Ops\", Word.InsertLocation.before);\nawait context.sync();\n\n// Restrict the \"not bold\" formatting to just the literal text we inserted\n// (not the paragraph mark) so the list-style/paragraph mark formatting\n// stays untouched, same as the sibling bullets in this list.\nconst devOpsRange = newPara.getRange().search(\"Dev Ops\", { matchCase: true });\ndevOpsRange.load(\"items\");\nawait context.sync();\ndevOpsRange.items[0].font.bold = false;\n\n// Append the bold \" - Moderate\" suffix as its own run.\nconst endRange = newPara.getRange(Word.RangeLocation.end);\nconst moderateRun = endRange.insertText(\" - Moderate\", Word.InsertLocation.end);\nmoderateRun.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Add a new \"Dev Ops - Moderate\" bullet to the Technical-Skills list, right\n# before the \"OOMD (Object Oriented Modeling and Design) - Basics\" item.\n# \"Dev Ops\" is regular weight, \" - Moderate\" is bold - matching the style\n# already used by every other skill line in this list.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its (stable) text rather than a hard-coded\n# index so the script keeps working even if earlier content shifts.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"OOMD (Object Oriented Modeling and Design)\")\nif (-not $found) {\n    throw 'Could not find the \"OOMD (Object Oriented Modeling and Design)\" paragraph.'\n}\n$anchorPara = $searchRange.Paragraphs(1)\n$anchorIndex = $anchorPara.Index\n\n# Insert a brand-new paragraph right before it; Word copies the anchor's\n# paragraph formatting (numbering, indent, fonts, size) automatically.\n$anchorPara.Range.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs($anchorIndex)\n$newRange = $newPara.Range\n$paraStart = $newRange.Start\n\n$devOpsText = \"Dev Ops\"\n$moderateText = \" - Moderate\"\n\n# Type the non-bold text first.\n$newRange.Text = $devOpsText\n\n# Restrict the \"not bold\" formatting to just the literal text we inserted\n# (not the paragraph mark) so the list-style/paragraph mark formatting\n# stays untouched, same as the sibling bullets in this list.\n$devOpsRange = $d.Range($paraStart, $paraStart + $devOpsText.Length)\n$devOpsRange.Font.Bold = $false\n\n# Append the bold \" - Moderate\" suffix as its own run.\n$insertionPoint = $d.Range($paraStart + $devOpsText.Length, $paraStart + $devOpsText.Length)\n$insertionPoint.InsertAfter($moderateText)\n$insertionPoint.Font.Bold = $true\n"}
